$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "1.0-Huevos (unidad),3.0-Leche (litros),2.0-Harina  (kg),1.0-Vainilla (ml),"
$ws.Range("C4").Value = "2.0-Huevos (unidad),5.0-Harina  (kg),1.0-Vainilla (ml),"
$ws.Range("C6").Value = "5.0-Crema (litros),5.0-Huevos (unidad),4.0-Harina  (kg),2.0-Limon (unidad),"
$ws.Range("C7").Value = "2.0-Huevos (unidad),0.2-Leche (litros),0.3-Harina  (kg),0.1-Vainilla (ml),"
